$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update the TCID value in A2 from "Ipa001" to "IPA001"
$ws.Range("A2").Value = "IPA001"
